$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (bottom-to-top so row numbers of not-yet-deleted rows stay stable):
#   17 - "cannot be assigned to a (non-viral) class" (to be replaced by a new combined "Unclassified" row)
#   20 - "unclassified"
#   21 - "Viruses"
#   19 - Eukaryota / Chytridiomycota / Neocallimastigomycetes
#   16 - Bacteria / Verrucomicrobia / Verrucomicrobiae
#   12 - Bacteria / Firmicutes / Erysipelotrichia
#    4 - Bacteria / Bacteroidetes / Chitinophagia
$ws.Rows(21).EntireRow.Delete() | Out-Null
$ws.Rows(20).EntireRow.Delete() | Out-Null
$ws.Rows(19).EntireRow.Delete() | Out-Null
$ws.Rows(17).EntireRow.Delete() | Out-Null
$ws.Rows(16).EntireRow.Delete() | Out-Null
$ws.Rows(12).EntireRow.Delete() | Out-Null
$ws.Rows(4).EntireRow.Delete() | Out-Null

# Add the new combined "Unclassified" row at the bottom (row 15)
$ws.Cells.Item(15, 1).Value = "Unclassified"
$ws.Cells.Item(15, 2).Value = "NA"
$ws.Cells.Item(15, 3).Value = "NA"
$ws.Cells.Item(15, 4).Value = "50.627 ± 2.190"
$ws.Cells.Item(15, 5).Value = "26.939 ± 6.526"
$ws.Cells.Item(15, 6).Value = "54.081 ± 1.963"

$ws.Range("A15").NumberFormat = "@"
$ws.Range("F15").NumberFormat = "0.00"

$ws.Range("F16").Select() | Out-Null
